$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. On the PRODOTTO worksheet, column C (TIPO_PRODOTTO_ID) currently
#    stores the numeric id of the product type (1..12). Replace every
#    value with the matching text code used on the TIPO_PRODOTTO sheet
#    ("T01".."T12"), keeping the row's other data untouched.
# ------------------------------------------------------------------
$wsProdotto = $wb.Worksheets.Item("PRODOTTO")

$lastRow = $wsProdotto.Cells.Item($wsProdotto.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $wsProdotto.Cells.Item($r, 3)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $code = "T" + "{0:D2}" -f [int]$val
        $cell.Value = $code
    }
}

# ------------------------------------------------------------------
# 2. The active/selected sheet moves from VENDITA to PRODOTTO
#    (tabSelected on the sheetView + workbook activeTab).
# ------------------------------------------------------------------
$wsProdotto.Activate()
